# "Inicio do Excel continua1"
# Populate the test-data sheet (Principal.xlsx) used by the TDD project:
# a header row (A1:C1) plus one row of sample credentials (B2:C2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("A1").Value = "preencheDadosUsuarioTesteCase"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "password"

# Data row
$ws.Range("B2").Value = "MafraMano"
$ws.Range("C2").Value = "4Jes"

# Column widths (character units, closest representable values)
$ws.Columns.Item(1).ColumnWidth = 39.5
$ws.Columns.Item(2).ColumnWidth = 12.333333333333334
$ws.Columns.Item(3).ColumnWidth = 12.5

# Select column A (mirrors the saved selection sqref="A1:A1048576")
$ws.Columns.Item(1).Select() | Out-Null
